# Optimizes 2 params for Pearson correlation
# - Recomputes B4 (measure for "one") from 3 to 2
# - Adds a new "two" row (row 5): measure goes from 9 to 4
# - Shifts the former "three" row down to row 6 with a new measure of 16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create the new row 6 first (re-using the existing "three" label) so the
# shared-string slot for "three" stays alive while we repurpose row 5 for "two".
$ws.Range("A6").Value = "three"
$ws.Range("B6").Value = 16

# Row 5 now becomes the "two" measure.
$ws.Range("A5").Value = "two"
$ws.Range("B5").Value = 4

# Row 4 ("one") measure is updated too.
$ws.Range("B4").Value = 2

# Move the active selection to B6, matching where the author left off editing.
$ws.Range("B6").Select() | Out-Null
